# Edit corresponding to commit "Sat, Jul 04, 2020  2:04:56 PM"
#
# 1) The table on slide 5 (the B1 "types of financial documents" table)
#    switches from table style {BF282AF5-99C2-4CB4-BE56-C7DEF1D8BBD6}
#    to {289482F3-AAEF-47BB-96DD-3A3A2DF1B04D}.
#
# 2) The deck's theme colour scheme (ppt/theme/theme2.xml, the theme
#    wired to the slide master / all slides) reverts from the
#    "Red Violet" / Integral palette back to the stock "Office" palette.

$p = $ppt.ActivePresentation

# --- 1) Table style -------------------------------------------------
for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if ($shape.HasTable) {
            $shape.Table.ApplyStyle("{289482F3-AAEF-47BB-96DD-3A3A2DF1B04D}")
        }
    }
}

# --- 2) Theme colours -------------------------------------------------
# MsoThemeColorSchemeIndex order: Dark1, Light1, Dark2, Light2,
# Accent1..Accent6, Hyperlink, FollowedHyperlink.
# Target values are the standard Office theme colours
# (dk1 000000, lt1 FFFFFF, dk2 44546A, lt2 E7E6E6, accent1 5B9BD5,
#  accent2 ED7D31, accent3 A5A5A5, accent4 FFC000, accent5 4472C4,
#  accent6 70AD47, hlink 0563C1, folHlink 954F72), expressed here as
# the decimal BGR values the PowerPoint RGB() macro would produce.
$officeThemeColors = @(
    0,         # dk1      000000
    16777215,  # lt1      FFFFFF
    6968388,   # dk2      44546A
    15132391,  # lt2      E7E6E6
    13998939,  # accent1  5B9BD5
    3243501,   # accent2  ED7D31
    10855845,  # accent3  A5A5A5
    49407,     # accent4  FFC000
    12874308,  # accent5  4472C4
    4697456,   # accent6  70AD47
    12673797,  # hlink    0563C1
    7491477    # folHlink 954F72
)

$themeColorScheme = $p.Slides.Item(1).ThemeColorScheme
for ($i = 1; $i -le $themeColorScheme.Count; $i++) {
    $themeColorScheme.Item($i).RGB = $officeThemeColors[$i - 1]
}
